# Update the bacterial-length figures on "Hoja1" (row 2) and refresh the
# selection/cursor position on both sheets, per the commit's data revision.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# Revised bacterial length figure + new central value / spread.
$ws1.Range("F2").Value = "1.4–2 µm long"
$ws1.Range("G2").Value = 1.7
$ws1.Range("H2").Value = 0.3

# Move the saved cursor position on Hoja2 (no longer a multi-cell selection).
$ws2.Range("B14").Select()

# Move the saved cursor position on Hoja1 last, so Hoja1 remains the
# workbook's active/visible tab on save.
$ws1.Range("F7").Select()
